# Update crypto price/volume data per the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.113.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.648.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.31%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +1.74%  "
$ws.Range("E10").Value = "  +6.88%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.125.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.972.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000147"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.642.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("E18").Value = "  +4.45%  "
$ws.Range("E19").Value = "  +3.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("E25").Value = "  +8.13%  "
$ws.Range("E26").Value = "  +4.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "556.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.09%  "
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("E33").Value = "  +5.14%  "
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.33"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "168.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.64%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "166.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.54%  "
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("E44").Value = "  +2.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0571"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("E48").Value = "  +2.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +13.19%  "
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.83%  "
Write-Host "Updated cryptos list with latest price/volume data."
